$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix double space typo in the "security guard" instructions cell (A3):
# "...responding as quickly as  possible..." -> "...responding as quickly as possible..."
$newText = "Imagine that you are a security guard watching for deviant activity at a business. Your job requires that you pay attention at all times, and respond quickly when something suspicious happens.  `n`nIn our lab we study attention and rapid responding, and in this experiment you'll be asked to play the role of the security guard. `n`nSpecifically, you will be attending to a number of items presented on the computer screen, and you'll be responding as quickly as possible when a target item appears by pressing the spacebar.  `n`n`nPress the spacebar to continue"

$ws.Range("A3").Value = $newText
